$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# K2 used to be a formula PRODUCT(6782,162253) (=1100399846); it becomes a
# plain numeric constant 1048576.
$ws.Range("K2").Value = 1048576

# H1's formula subtracts a different constant now.
$ws.Range("H1").Formula = "=(D1*K2)+(K2-1031670)"

# O49 held a leftover constant that is no longer needed; clear it so the
# used range shrinks back down to column K.
$ws.Range("O49").ClearContents()

# Re-enter the shared formula block so Excel regroups it as H3:H54 (H55
# keeps its own distinct, non-shared copy of the same formula).
$ws.Range("H3:H54").Formula = "=SUM(D3*`$K`$2,H2)"
$ws.Range("H55").Formula = "=SUM(D55*`$K`$2,H54)"

# Restore the saved selection (single cell K28, instead of the old H1:H55
# block selection).
$ws.Range("K28").Select()

$wb.Save()
